# Updates the cryptos table (rows 2-51) to the latest scrape: refreshed
# Price (column D) / Volume(1h) (column E) figures for every coin, plus a
# reordering of the last two rows (VeChain and Maker swap places, each
# carrying its own Coin/Link/Price/Volume values with it).
#
# Price/Volume are plain text in this sheet (e.g. "69.073.72", "  +2.61%  "),
# not numbers - that's how the source data is scraped/stored. Excel's COM
# layer auto-parses a bare numeric-looking string (e.g. "604.72") typed into
# Range.Value into a real number, which would both change the stored type
# and introduce binary float noise (604.72 -> 604.72000000000003) when
# that number is serialised back to XML. For the handful of D-column values
# that parse as plain numbers, force them to stay literal text the same way
# typing an apostrophe-prefixed entry in the Excel UI would: set the cell to
# Text format, assign the value, then clear the format again so the cell
# doesn't end up with a lingering NumberFormat override that wasn't there
# before.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Forces $value into $range as literal text even when it looks like a
    # number, without leaving a permanent number-format change behind.
    $ws.Range($range).NumberFormat = '@'
    $ws.Range($range).Value = $value
    $ws.Range($range).ClearFormats()
}


# Row 2
$ws.Range('D2').Value = '69.073.72'
$ws.Range('E2').Value = '  +2.61%  '

# Row 3
$ws.Range('D3').Value = '3.753.15'
$ws.Range('E3').Value = '  +2.24%  '

# Row 4
Set-TextValue 'D4' '0.999'
$ws.Range('E4').Value = '  -0.09%  '

# Row 5
Set-TextValue 'D5' '604.72'
$ws.Range('E5').Value = '  +1.62%  '

# Row 6
Set-TextValue 'D6' '169.43'
$ws.Range('E6').Value = '  +2.32%  '

# Row 7
$ws.Range('D7').Value = '3.747.15'
$ws.Range('E7').Value = '  +2.15%  '

# Row 8
$ws.Range('E8').Value = '  -0.02%  '

# Row 9
Set-TextValue 'D9' '0.537'
$ws.Range('E9').Value = '  +2.43%  '

# Row 10
Set-TextValue 'D10' '0.168'
$ws.Range('E10').Value = '  +6.26%  '

# Row 11
Set-TextValue 'D11' '6.37'
$ws.Range('E11').Value = '  +3.36%  '

# Row 12
Set-TextValue 'D12' '0.464'
$ws.Range('E12').Value = '  +0.78%  '

# Row 13
Set-TextValue 'D13' '38.44'
$ws.Range('E13').Value = '  +2.72%  '

# Row 14
$ws.Range('E14').Value = '  +4.34%  '

# Row 15
$ws.Range('D15').Value = '4.375.93'
$ws.Range('E15').Value = '  +2.21%  '

# Row 16
$ws.Range('D16').Value = '3.749.91'
$ws.Range('E16').Value = '  +2.18%  '

# Row 17
$ws.Range('D17').Value = '69.063.86'
$ws.Range('E17').Value = '  +2.50%  '

# Row 18
Set-TextValue 'D18' '7.30'
$ws.Range('E18').Value = '  +2.21%  '

# Row 19
$ws.Range('E19').Value = '  +0.08%  '

# Row 20
Set-TextValue 'D20' '17.14'
$ws.Range('E20').Value = '  -3.03%  '

# Row 21
Set-TextValue 'D21' '10.93'
$ws.Range('E21').Value = '  +20.09%  '

# Row 22
Set-TextValue 'D22' '495.72'
$ws.Range('E22').Value = '  +0.96%  '

# Row 23
Set-TextValue 'D23' '0.730'
$ws.Range('E23').Value = '  +1.67%  '

# Row 24
$ws.Range('E24').Value = '  +14.54%  '

# Row 25
Set-TextValue 'D25' '85.45'
$ws.Range('E25').Value = '  +0.15%  '

# Row 26
Set-TextValue 'D26' '2.34'
$ws.Range('E26').Value = '  +2.01%  '

# Row 27
Set-TextValue 'D27' '12.39'
$ws.Range('E27').Value = '  +2.23%  '

# Row 28
Set-TextValue 'D28' '10.36'
$ws.Range('E28').Value = '  +4.11%  '

# Row 29
$ws.Range('E29').Value = '  -0.09%  '

# Row 30
$ws.Range('E30').Value = '  +8.40%  '

# Row 31
$ws.Range('E31').Value = '  +2.69%  '

# Row 32
Set-TextValue 'D32' '7.99'
$ws.Range('E32').Value = '  +4.75%  '

# Row 33
Set-TextValue 'D33' '32.03'
$ws.Range('E33').Value = '  +2.15%  '

# Row 34
$ws.Range('D34').Value = '3.895.90'
$ws.Range('E34').Value = '  +2.29%  '

# Row 35
$ws.Range('E35').Value = '  +1.56%  '

# Row 36
$ws.Range('D36').Value = '3.683.74'
$ws.Range('E36').Value = '  +2.09%  '

# Row 37
$ws.Range('E37').Value = '  +0.03%  '

# Row 38
$ws.Range('E38').Value = '  +2.65%  '

# Row 39
Set-TextValue 'D39' '5.90'
$ws.Range('E39').Value = '  +2.77%  '

# Row 40
$ws.Range('E40').Value = '  +2.38%  '

# Row 41
Set-TextValue 'D41' '0.325'
$ws.Range('E41').Value = '  +1.11%  '

# Row 42
Set-TextValue 'D42' '3.06'
$ws.Range('E42').Value = '  +10.78%  '

# Row 43
Set-TextValue 'D43' '439.07'
$ws.Range('E43').Value = '  +1.57%  '

# Row 44
Set-TextValue 'D44' '48.83'
$ws.Range('E44').Value = '  +0.49%  '

# Row 45
Set-TextValue 'D45' '1.99'
$ws.Range('E45').Value = '  +3.19%  '

# Row 46
Set-TextValue 'D46' '8.48'
$ws.Range('E46').Value = '  +1.64%  '

# Row 47
$ws.Range('E47').Value = '  +0.00%  '

# Row 48
Set-TextValue 'D48' '40.47'
$ws.Range('E48').Value = '  +0.15%  '

# Row 49
Set-TextValue 'D49' '141.57'
$ws.Range('E49').Value = '  -0.78%  '

# Row 50
$ws.Range('B50').Value = 'Maker'
$ws.Range('C50').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D50').Value = '2.782.31'
$ws.Range('E50').Value = '  +1.24%  '

# Row 51
$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue 'D51' '0.0355'
$ws.Range('E51').Value = '  +2.64%  '
